$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86; this shifts the existing rows 86:164
# down to 87:165 (carrying their formatting with them).
$ws.Rows("86:86").Insert()

# Seed the new row 86 with the same data as the row below it (which is
# the original row 86, now shifted to row 87), then overwrite the
# columns that actually differ for the new record.
$vals = $ws.Range("A87:R87").Value2
$ws.Range("A86:R86").Value2 = $vals

$ws.Range("D86").Value2 = 44566
$ws.Range("J86").Value2 = 50
$ws.Range("K86").Value2 = 25000
$ws.Range("L86").Value2 = 25000
$ws.Range("M86").Value2 = 25000
$ws.Range("P86").Value2 = 1389
